$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 273.7846883333333
$ws.Cells.Item(2, 8).Value = 821.354065
$ws.Cells.Item(2, 9).Value = 0.8244324628389421
$ws.Cells.Item(2, 10).Value = 0.8244324628389421
$ws.Cells.Item(2, 13).Value = 211.2725676666666
$ws.Cells.Item(2, 14).Value = 633.8177029999999
$ws.Cells.Item(2, 15).Value = 0.6324644927232657
$ws.Cells.Item(2, 16).Value = 0.6324644927232657
$ws.Cells.Item(2, 17).Value = 57843.19409200141
$ws.Cells.Item(2, 18).Value = 520588.7468280126
$ws.Cells.Item(2, 19).Value = 0.521424259394024
$ws.Cells.Item(2, 20).Value = 0.521424259394024
$ws.Cells.Item(3, 7).Value = 273.7846883333333
$ws.Cells.Item(3, 8).Value = 821.354065
$ws.Cells.Item(3, 9).Value = 0.8244324628389421
$ws.Cells.Item(3, 10).Value = 0.8244324628389421
$ws.Cells.Item(3, 13).Value = 59.36675400000001
$ws.Cells.Item(3, 14).Value = 178.100262
$ws.Cells.Item(3, 15).Value = 0.1777200152765546
$ws.Cells.Item(3, 16).Value = 0.1777200152765546
$ws.Cells.Item(3, 17).Value = 16253.70824125167
$ws.Cells.Item(3, 18).Value = 146283.3741712651
$ws.Cells.Item(3, 19).Value = 0.1465181498902243
$ws.Cells.Item(3, 20).Value = 0.1465181498902243
$ws.Cells.Item(4, 7).Value = 273.7846883333333
$ws.Cells.Item(4, 8).Value = 821.354065
$ws.Cells.Item(4, 9).Value = 0.8244324628389421
$ws.Cells.Item(4, 10).Value = 0.8244324628389421
$ws.Cells.Item(4, 13).Value = 0.4593846666666666
$ws.Cells.Item(4, 14).Value = 1.378154
$ws.Cells.Item(4, 15).Value = 0.001375211620595172
$ws.Cells.Item(4, 16).Value = 0.001375211620595172
$ws.Cells.Item(4, 17).Value = 125.7724877884455
$ws.Cells.Item(4, 18).Value = 1131.95239009601
$ws.Cells.Item(4, 19).Value = 0.001133769103292011
$ws.Cells.Item(4, 20).Value = 0.001133769103292011
$ws.Cells.Item(5, 7).Value = 273.7846883333333
$ws.Cells.Item(5, 8).Value = 821.354065
$ws.Cells.Item(5, 9).Value = 0.8244324628389421
$ws.Cells.Item(5, 10).Value = 0.8244324628389421
$ws.Cells.Item(5, 13).Value = 62.94782133333333
$ws.Cells.Item(5, 14).Value = 188.843464
$ws.Cells.Item(5, 15).Value = 0.1884402803795846
$ws.Cells.Item(5, 16).Value = 0.1884402803795846
$ws.Cells.Item(5, 17).Value = 17234.14964500902
$ws.Cells.Item(5, 18).Value = 155107.3468050811
$ws.Cells.Item(5, 19).Value = 0.1553562844514017
$ws.Cells.Item(5, 20).Value = 0.1553562844514017
$ws.Cells.Item(6, 7).Value = 42.36118633333333
$ws.Cells.Item(6, 9).Value = 0.1275598624238963
$ws.Cells.Item(6, 10).Value = 0.1275598624238963
$ws.Cells.Item(6, 13).Value = 211.2725676666666
$ws.Cells.Item(6, 14).Value = 633.8177029999999
$ws.Cells.Item(6, 15).Value = 0.6324644927232657
$ws.Cells.Item(6, 16).Value = 0.6324644927232657
$ws.Cells.Item(6, 17).Value = 8949.756606049441
$ws.Cells.Item(6, 18).Value = 80547.80945444497
$ws.Cells.Item(6, 19).Value = 0.08067708367977913
$ws.Cells.Item(6, 20).Value = 0.08067708367977915
$ws.Cells.Item(7, 7).Value = 42.36118633333333
$ws.Cells.Item(7, 9).Value = 0.1275598624238963
$ws.Cells.Item(7, 10).Value = 0.1275598624238963
$ws.Cells.Item(7, 13).Value = 59.36675400000001
$ws.Cells.Item(7, 14).Value = 178.100262
$ws.Cells.Item(7, 15).Value = 0.1777200152765546
$ws.Cells.Item(7, 16).Value = 0.1777200152765546
$ws.Cells.Item(7, 17).Value = 2514.846128199162
$ws.Cells.Item(7, 18).Value = 22633.61515379246
$ws.Cells.Item(7, 19).Value = 0.02266994069865005
$ws.Cells.Item(7, 20).Value = 0.02266994069865005
$ws.Cells.Item(8, 7).Value = 42.36118633333333
$ws.Cells.Item(8, 9).Value = 0.1275598624238963
$ws.Cells.Item(8, 10).Value = 0.1275598624238963
$ws.Cells.Item(8, 13).Value = 0.4593846666666666
$ws.Cells.Item(8, 14).Value = 1.378154
$ws.Cells.Item(8, 15).Value = 0.001375211620595172
$ws.Cells.Item(8, 16).Value = 0.001375211620595172
$ws.Cells.Item(8, 17).Value = 19.46007946334288
$ws.Cells.Item(8, 18).Value = 175.140715170086
$ws.Cells.Item(8, 19).Value = 0.0001754218051268636
$ws.Cells.Item(8, 20).Value = 0.0001754218051268637
$ws.Cells.Item(9, 7).Value = 42.36118633333333
$ws.Cells.Item(9, 9).Value = 0.1275598624238963
$ws.Cells.Item(9, 10).Value = 0.1275598624238963
$ws.Cells.Item(9, 13).Value = 62.94782133333333
$ws.Cells.Item(9, 14).Value = 188.843464
$ws.Cells.Item(9, 15).Value = 0.1884402803795846
$ws.Cells.Item(9, 16).Value = 0.1884402803795846
$ws.Cells.Item(9, 17).Value = 2666.544388778708
$ws.Cells.Item(9, 18).Value = 23998.89949900837
$ws.Cells.Item(9, 19).Value = 0.02403741624034026
$ws.Cells.Item(9, 20).Value = 0.02403741624034026
$ws.Cells.Item(10, 7).Value = 1.581929
$ws.Cells.Item(10, 8).Value = 4.745787
$ws.Cells.Item(10, 9).Value = 0.004763573994753449
$ws.Cells.Item(10, 10).Value = 0.00476357399475345
$ws.Cells.Item(10, 13).Value = 211.2725676666666
$ws.Cells.Item(10, 14).Value = 633.8177029999999
$ws.Cells.Item(10, 15).Value = 0.6324644927232657
$ws.Cells.Item(10, 16).Value = 0.6324644927232657
$ws.Cells.Item(10, 17).Value = 334.2182016963623
$ws.Cells.Item(10, 18).Value = 3007.963815267261
$ws.Cells.Item(10, 19).Value = 0.00301279141014148
$ws.Cells.Item(10, 20).Value = 0.003012791410141481
$ws.Cells.Item(11, 7).Value = 1.581929
$ws.Cells.Item(11, 8).Value = 4.745787
$ws.Cells.Item(11, 9).Value = 0.004763573994753449
$ws.Cells.Item(11, 10).Value = 0.00476357399475345
$ws.Cells.Item(11, 13).Value = 59.36675400000001
$ws.Cells.Item(11, 14).Value = 178.100262
$ws.Cells.Item(11, 15).Value = 0.1777200152765546
$ws.Cells.Item(11, 16).Value = 0.1777200152765546
$ws.Cells.Item(11, 17).Value = 93.91398978846601
$ws.Cells.Item(11, 18).Value = 845.2259080961941
$ws.Cells.Item(11, 19).Value = 0.0008465824431185811
$ws.Cells.Item(11, 20).Value = 0.0008465824431185811
$ws.Cells.Item(12, 7).Value = 1.581929
$ws.Cells.Item(12, 8).Value = 4.745787
$ws.Cells.Item(12, 9).Value = 0.004763573994753449
$ws.Cells.Item(12, 10).Value = 0.00476357399475345
$ws.Cells.Item(12, 13).Value = 0.4593846666666666
$ws.Cells.Item(12, 14).Value = 1.378154
$ws.Cells.Item(12, 15).Value = 0.001375211620595172
$ws.Cells.Item(12, 16).Value = 0.001375211620595172
$ws.Cells.Item(12, 17).Value = 0.7267139263553332
$ws.Cells.Item(12, 18).Value = 6.540425337197999
$ws.Cells.Item(12, 19).Value = 0.000006550922313149909
$ws.Cells.Item(12, 20).Value = 0.00000655092231314991
$ws.Cells.Item(13, 7).Value = 1.581929
$ws.Cells.Item(13, 8).Value = 4.745787
$ws.Cells.Item(13, 9).Value = 0.004763573994753449
$ws.Cells.Item(13, 10).Value = 0.00476357399475345
$ws.Cells.Item(13, 13).Value = 62.94782133333333
$ws.Cells.Item(13, 14).Value = 188.843464
$ws.Cells.Item(13, 15).Value = 0.1884402803795846
$ws.Cells.Item(13, 16).Value = 0.1884402803795846
$ws.Cells.Item(13, 17).Value = 99.57898405401866
$ws.Cells.Item(13, 18).Value = 896.2108564861679
$ws.Cells.Item(13, 19).Value = 0.0008976492191802379
$ws.Cells.Item(13, 20).Value = 0.0008976492191802379
$ws.Cells.Item(14, 7).Value = 13.26179066666667
$ws.Cells.Item(14, 8).Value = 39.785372
$ws.Cells.Item(14, 9).Value = 0.03993448577249507
$ws.Cells.Item(14, 10).Value = 0.03993448577249507
$ws.Cells.Item(14, 13).Value = 211.2725676666666
$ws.Cells.Item(14, 14).Value = 633.8177029999999
$ws.Cells.Item(14, 15).Value = 0.6324644927232657
$ws.Cells.Item(14, 16).Value = 0.6324644927232657
$ws.Cells.Item(14, 17).Value = 2801.852566004502
$ws.Cells.Item(14, 18).Value = 25216.67309404051
$ws.Cells.Item(14, 19).Value = 0.02525714428626556
$ws.Cells.Item(14, 20).Value = 0.02525714428626556
$ws.Cells.Item(15, 7).Value = 13.26179066666667
$ws.Cells.Item(15, 8).Value = 39.785372
$ws.Cells.Item(15, 9).Value = 0.03993448577249507
$ws.Cells.Item(15, 10).Value = 0.03993448577249507
$ws.Cells.Item(15, 13).Value = 59.36675400000001
$ws.Cells.Item(15, 14).Value = 178.100262
$ws.Cells.Item(15, 15).Value = 0.1777200152765546
$ws.Cells.Item(15, 16).Value = 0.1777200152765546
$ws.Cells.Item(15, 17).Value = 787.3094641074962
$ws.Cells.Item(15, 18).Value = 7085.785176967465
$ws.Cells.Item(15, 19).Value = 0.007097157421549176
$ws.Cells.Item(15, 20).Value = 0.007097157421549174
$ws.Cells.Item(16, 7).Value = 13.26179066666667
$ws.Cells.Item(16, 8).Value = 39.785372
$ws.Cells.Item(16, 9).Value = 0.03993448577249507
$ws.Cells.Item(16, 10).Value = 0.03993448577249507
$ws.Cells.Item(16, 13).Value = 0.4593846666666666
$ws.Cells.Item(16, 14).Value = 1.378154
$ws.Cells.Item(16, 15).Value = 0.001375211620595172
$ws.Cells.Item(16, 16).Value = 0.001375211620595172
$ws.Cells.Item(16, 17).Value = 6.092263284809778
$ws.Cells.Item(16, 18).Value = 54.830369563288
$ws.Cells.Item(16, 19).Value = 0.00005491836889682778
$ws.Cells.Item(16, 20).Value = 0.00005491836889682778
$ws.Cells.Item(17, 7).Value = 13.26179066666667
$ws.Cells.Item(17, 8).Value = 39.785372
$ws.Cells.Item(17, 9).Value = 0.03993448577249507
$ws.Cells.Item(17, 10).Value = 0.03993448577249507
$ws.Cells.Item(17, 13).Value = 62.94782133333333
$ws.Cells.Item(17, 14).Value = 188.843464
$ws.Cells.Item(17, 15).Value = 0.1884402803795846
$ws.Cells.Item(17, 16).Value = 0.1884402803795846
$ws.Cells.Item(17, 17).Value = 834.800829445401
$ws.Cells.Item(17, 18).Value = 7513.207465008608
$ws.Cells.Item(17, 19).Value = 0.007525265695783503
$ws.Cells.Item(17, 20).Value = 0.007525265695783503
$ws.Cells.Item(18, 7).Value = 1.099085666666667
$ws.Cells.Item(18, 8).Value = 3.297257
$ws.Cells.Item(18, 9).Value = 0.003309614969913056
$ws.Cells.Item(18, 10).Value = 0.003309614969913057
$ws.Cells.Item(18, 13).Value = 211.2725676666666
$ws.Cells.Item(18, 14).Value = 633.8177029999999
$ws.Cells.Item(18, 15).Value = 0.6324644927232657
$ws.Cells.Item(18, 16).Value = 0.6324644927232657
$ws.Cells.Item(18, 17).Value = 232.2066508822967
$ws.Cells.Item(18, 18).Value = 2089.859857940671
$ws.Cells.Item(18, 19).Value = 0.002093213953055387
$ws.Cells.Item(18, 20).Value = 0.002093213953055388
$ws.Cells.Item(19, 7).Value = 1.099085666666667
$ws.Cells.Item(19, 8).Value = 3.297257
$ws.Cells.Item(19, 9).Value = 0.003309614969913056
$ws.Cells.Item(19, 10).Value = 0.003309614969913057
$ws.Cells.Item(19, 13).Value = 59.36675400000001
$ws.Cells.Item(19, 14).Value = 178.100262
$ws.Cells.Item(19, 15).Value = 0.1777200152765546
$ws.Cells.Item(19, 16).Value = 0.1777200152765546
$ws.Cells.Item(19, 17).Value = 65.249148397926
$ws.Cells.Item(19, 18).Value = 587.242335581334
$ws.Cells.Item(19, 19).Value = 0.0005881848230124622
$ws.Cells.Item(19, 20).Value = 0.0005881848230124622
$ws.Cells.Item(20, 7).Value = 1.099085666666667
$ws.Cells.Item(20, 8).Value = 3.297257
$ws.Cells.Item(20, 9).Value = 0.003309614969913056
$ws.Cells.Item(20, 10).Value = 0.003309614969913057
$ws.Cells.Item(20, 13).Value = 0.4593846666666666
$ws.Cells.Item(20, 14).Value = 1.378154
$ws.Cells.Item(20, 15).Value = 0.001375211620595172
$ws.Cells.Item(20, 16).Value = 0.001375211620595172
$ws.Cells.Item(20, 17).Value = 0.5049031026197777
$ws.Cells.Item(20, 18).Value = 4.544127923577999
$ws.Cells.Item(20, 19).Value = 0.000004551420966320176
$ws.Cells.Item(20, 20).Value = 0.000004551420966320177
$ws.Cells.Item(21, 7).Value = 1.099085666666667
$ws.Cells.Item(21, 8).Value = 3.297257
$ws.Cells.Item(21, 9).Value = 0.003309614969913056
$ws.Cells.Item(21, 10).Value = 0.003309614969913057
$ws.Cells.Item(21, 13).Value = 62.94782133333333
$ws.Cells.Item(21, 14).Value = 188.843464
$ws.Cells.Item(21, 15).Value = 0.1884402803795846
$ws.Cells.Item(21, 16).Value = 0.1884402803795846
$ws.Cells.Item(21, 17).Value = 69.18504817536088
$ws.Cells.Item(21, 18).Value = 622.665433578248
$ws.Cells.Item(21, 19).Value = 0.0006236647728788869
$ws.Cells.Item(21, 20).Value = 0.0006236647728788869
